$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 101 (shifts the old row 101..224 down to 102..225)
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row with the new price-record data
$ws.Cells.Item(101, 1).Value = 4
$ws.Cells.Item(101, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(101, 3).Value = "Los Lagos"
$ws.Cells.Item(101, 4).Value = 44546
$ws.Cells.Item(101, 5).Value = 10
$ws.Cells.Item(101, 6).Value = 100114014
$ws.Cells.Item(101, 7).Value = "Betarraga"
$ws.Cells.Item(101, 8).Value = "Sin especificar"
$ws.Cells.Item(101, 9).Value = "Primera"
$ws.Cells.Item(101, 10).Value = 500
$ws.Cells.Item(101, 11).Value = 800
$ws.Cells.Item(101, 12).Value = 1000
$ws.Cells.Item(101, 13).Value = 900
$ws.Cells.Item(101, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(101, 15).Value = "Región del Maule"
$ws.Cells.Item(101, 16).Value = 180
$ws.Cells.Item(101, 17).Value = 5
$ws.Cells.Item(101, 18).Value = "Hortaliza"
